$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the legend labels in row 1 to the "large-" prefixed dataset names
$ws.Range("B1").Value = "large-KDD99"
$ws.Range("C1").Value = "large-CoverType"
$ws.Range("D1").Value = "large-KDD98"

# Widen columns B and C to fit the new, longer labels
$ws.Columns.Item(2).ColumnWidth = 13.45
$ws.Columns.Item(3).ColumnWidth = 18.45

# Move the active cell selection
$ws.Range("D10").Select()
